$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("#0")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G39").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("G68").Value = 0

$ws = $wb.Worksheets.Item("#1")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("H2").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 0

$ws = $wb.Worksheets.Item("#2")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("F25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 0

$ws = $wb.Worksheets.Item("#3")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("F5").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("H18").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("G62").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("H71").Value = 0

$ws = $wb.Worksheets.Item("#4")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("H6").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("H49").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("F51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("G64").Value = 0
$ws.Range("F65").Value = 0

$ws = $wb.Worksheets.Item("#5")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("G24").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("F66").Value = 0
$ws.Range("H66").Value = 0

$ws = $wb.Worksheets.Item("#6")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("G3").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("G49").Value = 0
$ws.Range("G51").Value = 0
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H59").Value = 0
$ws.Range("F68").Value = 0
$ws.Range("F70").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("H71").Value = 0

$ws = $wb.Worksheets.Item("#7")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("H20").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("F35").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H42").Value = 0
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("F63").Value = 0
$ws.Range("G70").Value = 0

$ws = $wb.Worksheets.Item("#8")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("H14").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("H54").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H63").Value = 0

$ws = $wb.Worksheets.Item("#9")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("F7").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H48").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("H53").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("G67").Value = 0

$ws = $wb.Worksheets.Item("#10")
$ws.Range("I1").Value = "Motif #4"
$ws.Range("F9").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("H32").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("H40").Value = 0
$ws.Range("F52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("F54").Value = 0
$ws.Range("H59").Value = 0
$ws.Range("F64").Value = 0
$ws.Range("H64").Value = 0
$ws.Range("F69").Value = 0
$ws.Range("H69").Value = 0
